# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer "last updated" timestamp ---
$ws.Range("A1").Value2 = "Datos actualizados a 14 de Agosto de 2020 a las 06:24"

# --- Peru (row 9): Recuperados / Muertes updated ---
$ws.Range("E9").Value2 = 134342
$ws.Range("H9").Value2 = 25648

# --- Espana (row 13): Casos totales updated ---
$ws.Range("B13").Value2 = 355856

# --- Honduras (row 51): full row refresh, no re-rank ---
$ws.Range("B51").Value2 = 49042
$ws.Range("C51").Value2 = 385
$ws.Range("D51").Value2 = 7032
$ws.Range("E51").Value2 = 40468
$ws.Range("G51").Value2 = 9
$ws.Range("H51").Value2 = 1542

# --- Australia overtakes Austria (rows 71-72 swap) ---
$ws.Range("A71").Value2 = "Australia"
$ws.Range("B71").Value2 = 22741
$ws.Range("C71").Value2 = 383
$ws.Range("D71").Value2 = 13001
$ws.Range("E71").Value2 = 9365
$ws.Range("G71").Value2 = 14
$ws.Range("H71").Value2 = 375

$ws.Range("A72").Value2 = "Austria"
$ws.Range("B72").Value2 = 22594
$ws.Range("D72").Value2 = 20346
$ws.Range("E72").Value2 = 1523
$ws.Range("H72").Value2 = 725

# --- Libia overtakes Albania & Mauritania (rows 99-101 shift) ---
$ws.Range("A99").Value2 = "Libia"
$ws.Range("B99").Value2 = 7050
$ws.Range("D99").Value2 = 816
$ws.Range("E99").Value2 = 6099
$ws.Range("H99").Value2 = 135

$ws.Range("A100").Value2 = "Albania"
$ws.Range("B100").Value2 = 6971
$ws.Range("D100").Value2 = 3616
$ws.Range("E100").Value2 = 3142
$ws.Range("H100").Value2 = 213

$ws.Range("A101").Value2 = "Mauritania"
$ws.Range("B101").Value2 = 6653
$ws.Range("D101").Value2 = 5843
$ws.Range("E101").Value2 = 653
$ws.Range("H101").Value2 = 157

# --- Antigua y Barbuda overtakes San Martin (Parte Francesa) (rows 193-194 swap) ---
$ws.Range("A193").Value2 = "Antigua y Barbuda"
$ws.Range("B193").Value2 = 93
$ws.Range("C193").Value2 = 1
$ws.Range("D193").Value2 = 83
$ws.Range("E193").Value2 = 7
$ws.Range("H193").Value2 = 3

$ws.Range("A194").Value2 = "San Martin (Parte Francesa)"
$ws.Range("D194").Value2 = 45
$ws.Range("E194").Value2 = 43
$ws.Range("H194").Value2 = 4

# --- Islas Malvinas overtakes Montserrat (rows 213-214 swap) ---
$ws.Range("A213").Value2 = "Islas Malvinas"
$ws.Range("D213").Value2 = 13
$ws.Range("H213").Value2 = 0

$ws.Range("A214").Value2 = "Montserrat"
$ws.Range("D214").Value2 = 12
$ws.Range("H214").Value2 = 1
